# Common: Build edit works
# Adds translation rows for the "build" edit feature into the Import sheet,
# mirroring the existing "mixture" edit translation keys.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Extend formatting (style) of the last existing data row down across the
# new rows so the new cells pick up the same "import" cell style (s="1").
$ws.Range("A509:C509").Copy()
$ws.Range("A510:C519").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$rows = @(
    @{ Row = 510; Label = "lab.build.index.title";            Translation = "Detail buildu" },
    @{ Row = 511; Label = "lab.build.index.preview.title";     Translation = "Detail buildu" },
    @{ Row = 512; Label = "lab.build.index.preview.subtitle";  Translation = "Zde můžete spravovat vybraný build." },
    @{ Row = 513; Label = "lab.build.button.edit";             Translation = "Upravit" },
    @{ Row = 514; Label = "lab.build.preview.name";            Translation = "Jméno" },
    @{ Row = 515; Label = "lab.build.edit.title";              Translation = "Editace buildu" },
    @{ Row = 516; Label = "lab.build.edit.subtitle";           Translation = "Místo, kde je možné upravit build." },
    @{ Row = 517; Label = "lab.build.update.submit";           Translation = "Aktualizovat" },
    @{ Row = 518; Label = "lab.build.update.success";          Translation = "Build [{{data.name}}] byl aktualizován." },
    @{ Row = 519; Label = "lab.build.link.button";             Translation = "Zpět" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = "cs"
    $ws.Cells.Item($row, 2).Value = $r.Label
    $ws.Cells.Item($row, 3).Value = $r.Translation
}

# Restore the active selection to reflect where the user ended up editing.
$ws.Range("B514").Select()
